$wb = $excel.ActiveWorkbook

# Sheet 1: Battery_Data
$ws1 = $wb.Worksheets.Item("Battery_Data")
$ws1.Range("B2").Value = 406638.746988
$ws1.Range("B3").Value = 219584.92337352
$ws1.Range("B4").Value = 4391.698467470401
$ws1.Range("B5").Value = 21683.1714275

# Sheet 2: Yearly BRC
$ws2 = $wb.Worksheets.Item("Yearly BRC")
$ws2.Range("B2").Value = 5850.092220579987
$ws2.Range("B3").Value = 6936.583990863766
$ws2.Range("B4").Value = 6933.451650374368
$ws2.Range("B5").Value = 6621.112361776992
$ws2.Range("B6").Value = 6047.308640889105
